$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new Saturday value for the week starting row 6 (G6 = Saturday hours)
$ws.Range("G6").Value = 8.25

# Recalculate so formula cells (I6, I19, etc.) pick up the new value
$excel.Calculate()

# Update the active selection to match the recorded cursor position
$ws.Range("J8").Select()
